# Calibration.xlsx — "calibrated with new camera with new height"
#
# The camera was remounted at a new height (900 mm -> 1340 mm) and the
# pixel/mm calibration table was re-measured, so every "mm Distance"
# reading (column B, rows 2-28) changes. The scatter chart on the sheet
# plots Sheet1!$B$2:$B$28, so it reflects the same new numbers once Excel
# recalculates. The footer note in E30 documents the new camera height.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Updated calibration readings (mm Distance) for the new camera height ---
$ws.Range("B2").Value  = 18
$ws.Range("B3").Value  = 22
$ws.Range("B4").Value  = 29
$ws.Range("B5").Value  = 36
$ws.Range("B6").Value  = 43
$ws.Range("B7").Value  = 50
$ws.Range("B8").Value  = 58
$ws.Range("B9").Value  = 65
$ws.Range("B10").Value = 72
$ws.Range("B11").Value = 79
$ws.Range("B12").Value = 85
$ws.Range("B13").Value = 92
$ws.Range("B14").Value = 99
$ws.Range("B15").Value = 107
$ws.Range("B16").Value = 114
$ws.Range("B17").Value = 122
$ws.Range("B18").Value = 130
$ws.Range("B19").Value = 137
$ws.Range("B20").Value = 144
$ws.Range("B21").Value = 151
$ws.Range("B22").Value = 159
$ws.Range("B23").Value = 193
$ws.Range("B24").Value = 227
$ws.Range("B25").Value = 264
$ws.Range("B26").Value = 298
$ws.Range("B27").Value = 333
$ws.Range("B28").Value = 366

# --- Footer note: new camera height ---
$ws.Range("E30").Value = "height to the camera = 1340 mm"

# --- Refresh the chart so it picks up the new calibration numbers ---
$chartObj = $ws.ChartObjects().Item(1)
$chartObj.Chart.Refresh()
$wb.RefreshAll()
$excel.CalculateFullRebuild()

# --- Leave the cursor on the note cell, matching the last-saved selection ---
$ws.Range("E30").Select()
